$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 665; everything from row 665 down shifts to 666..748
$ws.Rows.Item(665).Insert()

# Populate the newly inserted row 665 with the new data record
$ws.Cells.Item(665, 1).Value = 10
$ws.Cells.Item(665, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(665, 3).Value = "La Araucanía"
$ws.Cells.Item(665, 4).Value = 44946
$ws.Cells.Item(665, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(665, 5).Value = 9
$ws.Cells.Item(665, 6).Value = 100112045
$ws.Cells.Item(665, 7).Value = "Zapallo"
$ws.Cells.Item(665, 8).Value = "Camote"
$ws.Cells.Item(665, 9).Value = "1a (cosecha)"
$ws.Cells.Item(665, 10).Value = 550
$ws.Cells.Item(665, 11).Value = 700
$ws.Cells.Item(665, 12).Value = 700
$ws.Cells.Item(665, 13).Value = 700
$ws.Cells.Item(665, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(665, 15).Value = "Región del Maule"
$ws.Cells.Item(665, 16).Value = 700
$ws.Cells.Item(665, 17).Value = 1
$ws.Cells.Item(665, 18).Value = "Hortaliza"
